$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: move C3 (and the rest of the new C:E block) onto A2's style FIRST. ---
# This frees up the old "B3/C3" shared style (index 5) so that, once only B3 is left
# on it, changing B3's alignment mutates that style in place instead of forking a
# brand new cellXfs entry.
$ws.Range("A2").Copy()
$ws.Range("C1:E3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Step 2: B3 is now the sole remaining user of the old wrap-centered style, so
# switching it to left-aligned rewrites that style definition in place. ---
$ws.Range("B3").HorizontalAlignment = -4131   # xlLeft

# --- Step 3: header/url styles for column F, copied verbatim (no re-alignment
# needed, so these just reuse the existing style indexes). ---
$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("F2").PasteSpecial(-4122)

# F3 picks up B3's *already updated* left-aligned wrap style.
$ws.Range("B3").Copy()
$ws.Range("F3").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Header row values ---
$ws.Range("C1").Value = "echo"
$ws.Range("D1").Value = "echo"
$ws.Range("E1").Value = "echo"
$ws.Range("F1").Value = "callApi"

# --- Row 2 values ---
$ws.Range("C2").Value = '${name}'
$ws.Range("D2").Value = '${age}'
$ws.Range("E2").Value = '${job}'

# --- Row 3: update B3 JSON (drop the "tests" block) and populate F3 with the new JSON ---
$ws.Range("B3").Value = "{`n  ""method"": ""GET"",`n  ""headers"": {`n    ""X-Redmine-API-Key"": ""ed7449e623fc4b7bef6b0353cd59c5af652e2d66""`n  },`n  ""body"": {},`n  ""params"": """",`n  ""store"": [`n    {`n      ""name"": ""name"",`n      ""value"": ""John""`n    },`n    {`n      ""name"": ""age"",`n      ""responseKey"": ""total_count""`n    },`n    {`n      ""name"": ""job"",`n      ""value"": ""Developer""`n    }`n  ]`n}"

$ws.Range("F3").Value = "{`n  ""target"": ""https://task.hugang.io/time_entries.json"",`n  ""value"": {`n    ""method"": ""GET"",`n    ""headers"": {`n      ""X-Redmine-API-Key"": ""ed7449e623fc4b7bef6b0353cd59c5af652e2d66""`n    },`n    ""body"": {},`n    ""params"": """",`n    ""store"": [`n      {`n        ""name"": ""name"",`n        ""value"": ""John""`n      },`n      {`n        ""name"": ""age"",`n        ""responseKey"": ""total_count""`n      },`n      {`n        ""name"": ""job"",`n        ""value"": ""Developer""`n      }`n    ]`n  }`n}"

# --- Row 3 height grows to Excel's max (409.5) to fit the longer F3 text ---
$ws.Rows.Item(3).RowHeight = 409.5

# --- Column widths: C:E shrink back toward the default width, F widens for the new long text ---
$ws.Columns.Item(3).ColumnWidth = 8.5
$ws.Columns.Item(4).ColumnWidth = 8.5
$ws.Columns.Item(5).ColumnWidth = 8.5
$ws.Columns.Item(6).ColumnWidth = 58.8

# --- Selection moves to F3 ---
$ws.Range("F3").Select()

$wb.Save()
